# Natmi following Dr Hou advice
# LR-pairs_lrc2p/Osm-Osmr.xlsx : expand the Osm/Osmr ligand-receptor table
# from 3 rows (ECs as the only sending cluster) to 6 rows (ECs and sCs as
# sending clusters, each against ECs/FAPs/sCs target clusters), refreshing
# every statistic column (E..T) for each pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# followed by the 16 numeric metric columns (E..T).
$rows = @(
    @("ECs", "Osm", "Osmr", "ECs",
      3, 1, 21.90566233333334, 65.716987, 0.9982377076397728, 0.9982377076397729,
      3, 1, 25.45779733333333, 76.373392, 0.2480357710869278, 0.2480357710869278,
      557.6699121344338, 5019.029209209904, 0.2475986595424782, 0.2475986595424783),
    @("ECs", "Osm", "Osmr", "FAPs",
      3, 1, 21.90566233333334, 65.716987, 0.9982377076397728, 0.9982377076397729,
      3, 1, 52.87700633333333, 158.631019, 0.515181610841248, 0.5151816108412481,
      1158.305845935528, 10424.75261341975, 0.514273710224333, 0.5142737102243331),
    @("ECs", "Osm", "Osmr", "sCs",
      3, 1, 21.90566233333334, 65.716987, 0.9982377076397728, 0.9982377076397729,
      3, 1, 24.30280066666667, 72.90840200000001, 0.2367826180718242, 0.2367826180718242,
      532.3689451583083, 4791.320506424775, 0.2363653378729616, 0.2363653378729616),
    @("sCs", "Osm", "Osmr", "ECs",
      1, 0.3333333333333333, 0.03867233333333333, 0.116017, 0.001762292360227098, 0.001762292360227098,
      3, 1, 25.45779733333333, 76.373392, 0.2480357710869278, 0.2480357710869278,
      0.984512424407111, 8.860611819663999, 0.0004371115444495301, 0.0004371115444495302),
    @("sCs", "Osm", "Osmr", "FAPs",
      1, 0.3333333333333333, 0.03867233333333333, 0.116017, 0.001762292360227098, 0.001762292360227098,
      3, 1, 52.87700633333333, 158.631019, 0.515181610841248, 0.5151816108412481,
      2.044877214591444, 18.403894931323, 0.0009079006169150212, 0.0009079006169150216),
    @("sCs", "Osm", "Osmr", "sCs",
      1, 0.3333333333333333, 0.03867233333333333, 0.116017, 0.001762292360227098, 0.001762292360227098,
      3, 1, 24.30280066666667, 72.90840200000001, 0.2367826180718242, 0.2367826180718242,
      0.939846008314889, 8.458614074834001, 0.0004172801988625465, 0.0004172801988625465)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 0; $c -lt $data.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $data[$c]
    }
}
